$d = $word.ActiveDocument

# 1. Split "Version 2" into two runs: "Version " and "3" (so the
#    revision number becomes its own run, matching the target OOXML).
$verPara = $null
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -match "^Version 2\r?\f?$") { $verPara = $p; break }
}
if ($verPara -ne $null) {
    $verXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r><w:t xml:space="preserve">Version </w:t></w:r>' +
              '<w:r><w:t>3</w:t></w:r>' +
              '</w:p>'
    $verPara.Range.InsertXML($verXml)
}

# 2. Update the date field's cached text.
$d.Content.Find.Execute("5/26/21 12:31 PM", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11/15/21 10:01 AM", 2)

# 3. Remove the "Document Status" heading paragraph and the
#    "Approved by certification body" paragraph underneath it, while
#    keeping the heading's paragraph mark (its pPr/bookmark) attached
#    to the following "License" heading text. Find the three
#    paragraphs by their text and replace the whole span with a single
#    paragraph that carries the first paragraph's formatting/bookmark
#    and the third paragraph's run.
$statusPara = $null
$approvedPara = $null
$licensePara = $null
foreach ($p in @($d.Paragraphs)) {
    $t = $p.Range.Text
    if ($t -match "^Document Status\r?\f?$") { $statusPara = $p }
    elseif ($t -match "^Approved by certification body\r?\f?$") { $approvedPara = $p }
    elseif ($statusPara -ne $null -and $approvedPara -ne $null -and $licensePara -eq $null -and $t -match "^License\r?\f?$") { $licensePara = $p }
}

if ($statusPara -ne $null -and $approvedPara -ne $null -and $licensePara -ne $null) {
    $span = $d.Range($statusPara.Range.Start, $licensePara.Range.End)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="both"/></w:pPr>' +
           '<w:bookmarkStart w:id="400" w:name="_f5diy2ktdyyf" w:colFirst="0" w:colLast="0"/>' +
           '<w:bookmarkEnd w:id="400"/>' +
           '<w:r><w:t>License</w:t></w:r>' +
           '</w:p>'
    $span.InsertXML($xml)
}
